# Change to use offline signing for ManOps 2 New Sale and List contracts wip
#
# Insert a new "combo" bit-flag row above the current row 54 on the
# "List Entry Bits" sheet:
#   LE_SALE_CON_PICOS_FR_TRAN_OK_B =    1034;
#       // LE_SALE_CONTRACT_B | LE_HOLDS_PICOS_B | LE_FROM_TRANSFER_OK_B
#       // for the sale contract bit settings
#
# Inserting the row shifts every row from the old 54 downward by one,
# and Excel automatically re-points any formula that referenced the old
# B54 (e.g. the old row 57's "A29+A27+A31+B54" becomes
# "A29+A27+A31+B55" on the new row 58), exactly matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the existing row 54, pushing rows 54-62
# down to 55-63 and shifting their relative formula references.
$ws.Rows("54:54").Insert()

# B54: the new combo value = LE_SALE_CONTRACT_B | LE_HOLDS_PICOS_B | LE_FROM_TRANSFER_OK_B
$ws.Range("B54").Formula = "=A25+A27+A34"

# C54: the new description string, formatted like the other "combo" rows
# (the bold/quote-prefixed Lucida Console style used on row 53, C53).
$ws.Range("C54").Value = " LE_SALE_CON_PICOS_FR_TRAN_OK_B =    1034; // LE_SALE_CONTRACT_B | LE_HOLDS_PICOS_B | LE_FROM_TRANSFER_OK_B for the sale contract bit settings"
$ws.Range("C53").Copy()
$ws.Range("C54").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the sheet's selection to reflect the last-edited cell.
$ws.Range("B54").Select()
